# "Daily task" workbook update — log entries for 2025-01-17 .. 2025-01-28
# (rows 14-25 on the "python" sheet), plus tidy-up of the placeholder rows
# 8-12 that were sitting empty/pre-formatted ahead of the typed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Clean up rows 8-12: they were pre-formatted (date style + an unused
#    "applyAlignment" style on B) but B8,B9,B11,B12 were never actually
#    filled in, and B10 just needs its format stripped back to plain.
# ---------------------------------------------------------------------
$ws.Range("B8").Clear()
$ws.Range("B9").Clear()
$ws.Range("B11").Clear()
$ws.Range("B12").Clear()
$ws.Range("B10").ClearFormats()

# ---------------------------------------------------------------------
# 2) New rows of data. Values are entered in the exact order below so
#    that shared-string table slots line up with how the sheet was
#    actually typed (row 20's text was entered before row 19's).
# ---------------------------------------------------------------------
$ws.Range("A14").Value = 45674
$ws.Range("B14").Value = "dataiku recepes practice"
$ws.Range("D14").Value = "completed"

$ws.Range("A15").Value = 45675
$ws.Range("B15").Value = "leave"

$ws.Range("A16").Value = 45676
$ws.Range("B16").Value = "Sunday"

$ws.Range("A17").Value = 45677
$ws.Range("B17").Value = "Monday-client holiday"

$ws.Range("A18").Value = 45678
$ws.Range("B18").Value = "dataiku recepes group, join recepes"
$ws.Range("D18").Value = "completed"

$ws.Range("A20").Value = 45680
$ws.Range("B20").Value = "filter, sync, distinct recepes "
$ws.Range("D20").Value = "completed"

$ws.Range("A19").Value = 45679
$ws.Range("B19").Value = "python ,add recepes"
$ws.Range("D19").Value = "completed"

$ws.Range("A21").Value = 45681
$ws.Range("B21").Value = "split , prepare recepes"
$ws.Range("D21").Value = "completed"

$ws.Range("A22").Value = 45682
$ws.Range("B22").Value = "dataiku recepes practice"
$ws.Range("D22").Value = "completed"

$ws.Range("A23").Value = 45683
$ws.Range("B23").Value = "sunday"

$ws.Range("A24").Value = 38379
$ws.Range("B24").Value = "parameters configurations ,formatt recepe"

$ws.Range("A25").Value = 45685

# ---------------------------------------------------------------------
# 3) Match the date formatting already used in column A (numFmt 14,
#    "completed"-row style) for both the tidied rows and the new rows,
#    by copying the format from A2 instead of assigning a NumberFormat
#    string (which would mint a brand-new style entry).
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A8:A12").PasteSpecial(-4122)
$ws.Range("A14:A25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Leave the view on the newly-entered last row, scrolled so row 6 is
#    at the top of the window.
# ---------------------------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A25").Select()
$excel.ActiveWindow.ScrollRow = 6
